$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text cells (Coin name / Link / Volume%) that are safe to assign directly as strings.
$textUpdates = @{
    'E2' = '  +0.84%  '
    'E3' = '  +0.86%  '
    'E4' = '  -0.11%  '
    'E5' = '  +1.43%  '
    'E6' = '  -2.09%  '
    'E7' = '  +7.21%  '
    'E8' = '  -0.13%  '
    'E9' = '  +7.87%  '
    'E10' = '  +10.84%  '
    'E11' = '  +2.15%  '
    'E12' = '  +0.16%  '
    'B13' = 'Polkadot'
    'C13' = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    'E13' = '  +8.72%  '
    'B14' = 'WrappedliquidstakedEther2.0'
    'C14' = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    'E14' = '  +0.49%  '
    'E15' = '  +7.51%  '
    'E16' = '  +44.74%  '
    'E17' = '  -0.18%  '
    'E18' = '  +6.66%  '
    'E19' = '  +6.77%  '
    'E20' = '  +0.62%  '
    'E21' = '  +42.32%  '
    'E22' = '  +9.83%  '
    'E23' = '  +1.77%  '
    'E24' = '  +1.75%  '
    'E25' = '  +3.19%  '
    'E26' = '  +11.68%  '
    'E27' = '  +9.25%  '
    'E28' = '  -0.25%  '
    'E29' = '  +0.79%  '
    'E30' = '  -5.07%  '
    'E31' = '  +5.50%  '
    'E32' = '  +0.12%  '
    'E33' = '  -0.33%  '
    'E34' = '  -1.31%  '
    'E35' = '  +0.02%  '
    'E36' = '  +3.82%  '
    'E37' = '  +3.43%  '
    'E38' = '  -0.18%  '
    'E39' = '  +1.09%  '
    'E40' = '  +7.78%  '
    'E41' = '  -0.10%  '
    'E42' = '  -1.59%  '
    'E43' = '  +1.53%  '
    'E44' = '  +8.31%  '
    'E45' = '  +1.35%  '
    'E46' = '  +8.43%  '
    'E47' = '  +0.27%  '
    'E48' = '  +5.99%  '
    'E49' = '  +0.70%  '
    'E50' = '  +1.34%  '
    'E51' = '  +7.50%  '
}

foreach ($cell in $textUpdates.Keys) {
    $ws.Range($cell).Value = $textUpdates[$cell]
}

# Price cells (column D). Many of these look like plain numbers to Excel (e.g. "128.27")
# and would otherwise be auto-converted to a floating point number (losing trailing
# zeroes / introducing rounding artifacts). Force them to remain plain text, matching
# the original inlineStr cell type, then restore the default (unstyled) cell format so
# we do not leave a stray number format behind.
$priceUpdates = @{
    'D2' = '61.903.43'
    'D3' = '3.420.83'
    'D6' = '128.27'
    'D7' = '0.630'
    'D9' = '0.733'
    'D11' = '42.42'
    'D13' = '9.04'
    'D14' = '3.966.00'
    'D16' = '0.0000201'
    'D17' = '3.430.74'
    'D18' = '12.39'
    'D19' = '1.07'
    'D20' = '61.871.84'
    'D21' = '441.35'
    'D22' = '90.94'
    'D24' = '12.89'
    'D25' = '3.23'
    'D26' = '32.80'
    'D27' = '8.71'
    'D29' = '2.75'
    'D30' = '7.55'
    'D31' = '11.90'
    'D33' = '0.114'
    'D34' = '42.57'
    'D36' = '0.0497'
    'D37' = '53.05'
    'D38' = '0.999'
    'D39' = '3.38'
    'D42' = '0.312'
    'D43' = '141.36'
    'D44' = '4.25'
    'D47' = '16.51'
    'D48' = '22.43'
    'D49' = '3.770.55'
    'D50' = '2.121.45'
    'D51' = '2.07'
}

foreach ($cell in $priceUpdates.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $priceUpdates[$cell]
    $range.Style = "Normal"
}
